$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "a"

# Row 2
$ws.Range("A2").Value = "b"
$ws.Range("B2").Value = "'1"
$ws.Range("C2").Value = "'2"

# Row 3
$ws.Range("A3").Value = "c"

# Row 4
$ws.Range("A4").Value = "d"

# Row 5
$ws.Range("A5").Value = "'5"
$ws.Range("B5").Value = "'8"
$ws.Range("C5").Value = "'9"

# Row 6
$ws.Range("A6").Value = "'6"
